# The sheet is a daily-price log for "Acelga" at Femacal de La Calera.
# A new weekly record was inserted as row 149 (pushing every following
# row down by one, through the former row 232 which becomes row 233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 149; everything at/after 149 shifts
# down one row (old row 149 -> 150, ..., old row 232 -> 233).
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new record's data.
$ws.Range("A149").Value = 3
$ws.Range("B149").Value = "Femacal de La Calera"
$ws.Range("C149").Value = "Coquimbo"
$ws.Range("D149").Value = 44518
$ws.Range("E149").Value = 5
$ws.Range("F149").Value = 100112009
$ws.Range("G149").Value = "Acelga"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 310
$ws.Range("K149").Value = 2000
$ws.Range("L149").Value = 2200
$ws.Range("M149").Value = 2097
$ws.Range("N149").Value = '$/docena de atados (6 kilos)'
$ws.Range("O149").Value = "Provincia de Quillota"
$ws.Range("P149").Value = 350
$ws.Range("Q149").Value = 6
$ws.Range("R149").Value = "Hortaliza"
